# Update "想去人数" (column F) figures for the gh-pages data refresh.
# Source: commit "Update gh-pages to output generated at 456a3b4"
$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 87
$ws.Cells.Item(4, 6).Value = 1795
$ws.Cells.Item(5, 6).Value = 3351
$ws.Cells.Item(6, 6).Value = 1156
$ws.Cells.Item(7, 6).Value = 2259
$ws.Cells.Item(8, 6).Value = 2184
$ws.Cells.Item(9, 6).Value = 1147
$ws.Cells.Item(10, 6).Value = 621
$ws.Cells.Item(12, 6).Value = 1706
$ws.Cells.Item(14, 6).Value = 78
$ws.Cells.Item(15, 6).Value = 50
$ws.Cells.Item(18, 6).Value = 1617
$ws.Cells.Item(19, 6).Value = 278
$ws.Cells.Item(20, 6).Value = 1338
$ws.Cells.Item(21, 6).Value = 754
$ws.Cells.Item(22, 6).Value = 283
$ws.Cells.Item(23, 6).Value = 636
$ws.Cells.Item(24, 6).Value = 12432
$ws.Cells.Item(25, 6).Value = 12475
$ws.Cells.Item(26, 6).Value = 924
$ws.Cells.Item(27, 6).Value = 718
$ws.Cells.Item(29, 6).Value = 264
$ws.Cells.Item(30, 6).Value = 38
$ws.Cells.Item(31, 6).Value = 416
$ws.Cells.Item(32, 6).Value = 1941
$ws.Cells.Item(33, 6).Value = 8
$ws.Cells.Item(35, 6).Value = 216
$ws.Cells.Item(36, 6).Value = 632

# --- Sheet 2: 演出 ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(10, 6).Value = 57

# --- Sheet 3: 本地生活 ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 6).Value = 117

# --- Sheet 4: 全部类型 ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 87
$ws.Cells.Item(5, 6).Value = 1795
$ws.Cells.Item(6, 6).Value = 3351
$ws.Cells.Item(7, 6).Value = 1156
$ws.Cells.Item(8, 6).Value = 2259
$ws.Cells.Item(9, 6).Value = 2184
$ws.Cells.Item(10, 6).Value = 1147
$ws.Cells.Item(11, 6).Value = 621
$ws.Cells.Item(12, 6).Value = 117
$ws.Cells.Item(14, 6).Value = 1706
$ws.Cells.Item(17, 6).Value = 78
$ws.Cells.Item(18, 6).Value = 50
$ws.Cells.Item(23, 6).Value = 1617
$ws.Cells.Item(24, 6).Value = 278
$ws.Cells.Item(25, 6).Value = 1338
$ws.Cells.Item(26, 6).Value = 754
$ws.Cells.Item(27, 6).Value = 283
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(29, 6).Value = 636
$ws.Cells.Item(30, 6).Value = 12432
$ws.Cells.Item(31, 6).Value = 12475
$ws.Cells.Item(32, 6).Value = 924
$ws.Cells.Item(33, 6).Value = 718
$ws.Cells.Item(35, 6).Value = 264
$ws.Cells.Item(36, 6).Value = 38
$ws.Cells.Item(37, 6).Value = 416
$ws.Cells.Item(40, 6).Value = 1941
$ws.Cells.Item(41, 6).Value = 8
$ws.Cells.Item(45, 6).Value = 216
$ws.Cells.Item(46, 6).Value = 632
$ws.Cells.Item(47, 6).Value = 57

